# Apply the commit: update header label in B1 from "commentaire" to "nom"
# and move the active selection to B2 (as captured in the sheet view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "nom"

$ws.Range("B2").Select()
